# Remove 2nd TCXO as it is not needed
# This removes reference "X2" (and "C37" / "C53" from shared BOM rows that
# shared a capacitor designator with the 2nd TCXO circuit) and decrements
# the affected QTY values by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: C37 removed from Refs list, QTY 13 -> 12
$ws.Range("A4").Value = "C49;C39;C16;C18;C4;C32;C42;C38;C15;C17;C3;C31"
$ws.Range("D4").Value = 12

# Row 12: C53 removed from Refs list, QTY 3 -> 2
$ws.Range("A12").Value = "C34;C33"
$ws.Range("D12").Value = 2

# Row 19: X2 removed from Refs list, QTY 2 -> 1
$ws.Range("A19").Value = "X1"
$ws.Range("D19").Value = 1

# Update the selection to reflect where the edit was made
$ws.Range("A12").Select() | Out-Null
